# 10.b.1.1 workbook fix-up:
#   - The sheet title cells (A1 Kyrgyz, C1 English) were missing the ".1"
#     suffix ("10.b.1" -> "10.b.1.1") that the Russian title (B1) already
#     had. Re-write those two cells with the corrected text.
#   - Update the active selection to L8, matching the re-saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "10.b.1.1 Агымдардын түрлөрү жана алуучу өлкөлөр жана донор-өлкөлөр боюнча бөлунүшүндөгү  өнүктүрүү максатында ресурстар агымынын жалпы көлөмү"
$ws.Range("C1").Value = "10.b.1.1 Total resource flows for development, by recipient and donor countries and type of flow (e.g. official development assistance, foreign direct investment and other flows)"

$ws.Range("L8").Select()
